$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.047125
$ws.Range("H2").Value = 0.141375
$ws.Range("I2").Value = 0.1108387998127795
$ws.Range("J2").Value = 0.1108387998127795
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 19.11595033333333
$ws.Range("N2").Value = 57.347851
$ws.Range("O2").Value = 0.6851940154453416
$ws.Range("P2").Value = 0.6851940154453418
$ws.Range("Q2").Value = 0.9008391594583334
$ws.Range("R2").Value = 8.107552435125
$ws.Range("S2").Value = 0.07594608231086074
$ws.Range("T2").Value = 0.07594608231086075

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.047125
$ws.Range("H3").Value = 0.141375
$ws.Range("I3").Value = 0.1108387998127795
$ws.Range("J3").Value = 0.1108387998127795
$ws.Range("M3").Value = 4.865208333333334
$ws.Range("O3").Value = 0.1743890089566637
$ws.Range("P3").Value = 0.1743890089566637
$ws.Range("Q3").Value = 0.2292729427083334
$ws.Range("R3").Value = 2.063456484375
$ws.Range("S3").Value = 0.01932906845329665
$ws.Range("T3").Value = 0.01932906845329666

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.047125
$ws.Range("H4").Value = 0.141375
$ws.Range("I4").Value = 0.1108387998127795
$ws.Range("J4").Value = 0.1108387998127795
$ws.Range("M4").Value = 3.917436333333333
$ws.Range("N4").Value = 11.752309
$ws.Range("O4").Value = 0.1404169755979945
$ws.Range("P4").Value = 0.1404169755979946
$ws.Range("Q4").Value = 0.1846091872083333
$ws.Range("R4").Value = 1.661482684875
$ws.Range("S4").Value = 0.01556364904862205
$ws.Range("T4").Value = 0.01556364904862206

# Row 5
$ws.Range("I5").Value = 0.5286385506557817
$ws.Range("J5").Value = 0.5286385506557816
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 19.11595033333333
$ws.Range("N5").Value = 57.347851
$ws.Range("O5").Value = 0.6851940154453416
$ws.Range("P5").Value = 0.6851940154453418
$ws.Range("Q5").Value = 4.296494624936556
$ws.Range("R5").Value = 38.668451624429
$ws.Range("S5").Value = 0.3622199712430407
$ws.Range("T5").Value = 0.3622199712430407

# Row 6
$ws.Range("I6").Value = 0.5286385506557817
$ws.Range("J6").Value = 0.5286385506557816
$ws.Range("M6").Value = 4.865208333333334
$ws.Range("O6").Value = 0.1743890089566637
$ws.Range("P6").Value = 0.1743890089566637
$ws.Range("S6").Value = 0.09218875294514883
$ws.Range("T6").Value = 0.09218875294514883

# Row 7
$ws.Range("I7").Value = 0.5286385506557817
$ws.Range("J7").Value = 0.5286385506557816
$ws.Range("M7").Value = 3.917436333333333
$ws.Range("N7").Value = 11.752309
$ws.Range("O7").Value = 0.1404169755979945
$ws.Range("P7").Value = 0.1404169755979946
$ws.Range("Q7").Value = 0.880481684467889
$ws.Range("R7").Value = 7.924335160211001
$ws.Range("S7").Value = 0.0742298264675921
$ws.Range("T7").Value = 0.0742298264675921

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.1532823333333333
$ws.Range("H8").Value = 0.459847
$ws.Range("I8").Value = 0.360522649531439
$ws.Range("J8").Value = 0.360522649531439
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 19.11595033333333
$ws.Range("N8").Value = 57.347851
$ws.Range("O8").Value = 0.6851940154453416
$ws.Range("P8").Value = 0.6851940154453418
$ws.Range("Q8").Value = 2.930137470977444
$ws.Range("R8").Value = 26.371237238797
$ws.Range("S8").Value = 0.2470279618914403
$ws.Range("T8").Value = 0.2470279618914403

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.1532823333333333
$ws.Range("H9").Value = 0.459847
$ws.Range("I9").Value = 0.360522649531439
$ws.Range("J9").Value = 0.360522649531439
$ws.Range("M9").Value = 4.865208333333334
$ws.Range("O9").Value = 0.1743890089566637
$ws.Range("P9").Value = 0.1743890089566637
$ws.Range("Q9").Value = 0.7457504854861112
$ws.Range("R9").Value = 6.711754369375
$ws.Range("S9").Value = 0.06287118755821824
$ws.Range("T9").Value = 0.06287118755821826

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.1532823333333333
$ws.Range("H10").Value = 0.459847
$ws.Range("I10").Value = 0.360522649531439
$ws.Range("J10").Value = 0.360522649531439
$ws.Range("M10").Value = 3.917436333333333
$ws.Range("N10").Value = 11.752309
$ws.Range("O10").Value = 0.1404169755979945
$ws.Range("P10").Value = 0.1404169755979946
$ws.Range("Q10").Value = 0.6004737818581111
$ws.Range("R10").Value = 5.404264036722999
$ws.Range("S10").Value = 0.05062350008178041
$ws.Range("T10").Value = 0.05062350008178042
